# Apply weekly update: rotate the Fecha/Volumen/Precio/Unidad/Origen/Kg data
# across rows 2-20 (row 6 untouched) following the permutation observed in
# the target diff. For each destination row we write the values that
# originally lived in the corresponding source row (captured below from the
# "before" state), so the move is applied consistently in one pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original ("before") values keyed by row number, for columns:
# D=Fecha, M=Volumen, N=Precio minimo, O=Precio maximo, P=Precio promedio ponderado,
# Q=Unidad de comercializacion, R=Origen, S=Precio $/Kg, T=Kg/unidad
$original = @{
  2  = @{ D = 44511; M = 45;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
  3  = @{ D = 44511; M = 45;  N = 3200;  O = 3200;  P = 3200;  Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 320;  T = 10 }
  4  = @{ D = 44859; M = 30;  N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 4000; T = 5  }
  5  = @{ D = 44519; M = 30;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
  7  = @{ D = 44889; M = 50;  N = 30000; O = 30000; P = 30000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 3000; T = 10 }
  8  = @{ D = 44503; M = 50;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
  9  = @{ D = 44515; M = 80;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 }
  10 = @{ D = 44488; M = 100; N = 12000; O = 12000; P = 12000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";               S = 2400; T = 5  }
  11 = @{ D = 44902; M = 90;  N = 25000; O = 25000; P = 25000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2500; T = 10 }
  12 = @{ D = 44483; M = 35;  N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 2000; T = 5  }
  13 = @{ D = 44466; M = 80;  N = 11000; O = 11000; P = 11000; Q = "`$/bandeja 5 kilos";  R = "La Ligua";               S = 2200; T = 5  }
  14 = @{ D = 44874; M = 40;  N = 25000; O = 25000; P = 25000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2500; T = 10 }
  15 = @{ D = 44858; M = 90;  N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 4000; T = 5  }
  16 = @{ D = 44496; M = 55;  N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2800; T = 10 }
  17 = @{ D = 44879; M = 25;  N = 30000; O = 30000; P = 30000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 3000; T = 10 }
  18 = @{ D = 44868; M = 30;  N = 14000; O = 14000; P = 14000; Q = "`$/bandeja 5 kilos";  R = "Provincia de Quillota";  S = 2800; T = 5  }
  19 = @{ D = 44921; M = 55;  N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 7 kilos";  R = "Provincia de Quillota";  S = 2143; T = 7  }
  20 = @{ D = 44901; M = 40;  N = 25000; O = 25000; P = 25000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Quillota";  S = 2500; T = 10 }
}

# Destination row -> source row (data that ends up living in destination row)
$mapping = @{
  2  = 5
  3  = 10
  4  = 18
  5  = 20
  7  = 16
  8  = 13
  9  = 7
  10 = 12
  11 = 17
  12 = 19
  13 = 2
  14 = 3
  15 = 4
  16 = 9
  17 = 15
  18 = 8
  19 = 14
  20 = 11
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $vals = $original[$srcRow]

    $ws.Cells.Item($destRow, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($destRow, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($destRow, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($destRow, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($destRow, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($destRow, 17).Value = $vals.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($destRow, 18).Value = $vals.R   # R - Origen
    $ws.Cells.Item($destRow, 19).Value = $vals.S   # S - Precio $/Kg
    $ws.Cells.Item($destRow, 20).Value = $vals.T   # T - Kg / unidad
}
